$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("campusTest")

$ws.Range("A78").Value = "Monitorizing assignments in Assign Feature"
$ws.Range("B78").Value = "FAILED"
$ws.Range("C78").Value = "chrome"

$ws.Range("A79").Value = "Monitorizing assignments in Assign Feature"
$ws.Range("B79").Value = "PASSED"
$ws.Range("C79").Value = "chrome"
